# Applies the "Small map changes, more diss work" edit to the Methodology
# section, renames the "Implementation & Results & Analysis" heading to
# "Critical Reflection", and adds a lastRenderedPageBreak marker before the
# "What have I accomplished?" paragraph.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Expand the five brief Methodology bullet paragraphs into the new,
#    much more detailed outline (13 paragraphs, some indented/tabbed).
# ---------------------------------------------------------------------------

$rStart = $d.Content
$rStart.Find.Execute("Planning, structure, and initial plans for the project.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $rStart.Start

$rEnd = $d.Content
$rEnd.Find.Execute("Landscape generation technology & real-time mesh deformation, 3D representation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $rEnd.End

$outlineBody = ""
$outlineBody += "<w:p><w:r><w:t>Planning</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>Class Structure</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>Noise &amp; Heightmap Generation</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>Soilmap Generation</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>OpenGL visuals</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:tab/><w:t>Mesh deformation</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:tab/><w:t>Visualization of fluid behaviour- pools</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:tab/><w:t>Visualisation of fluid behaviour- particles</w:t></w:r><w:r><w:t xml:space=`"preserve`"> &amp; foliage</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>Drop &amp; fluid as a particle</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr><w:r><w:t>Descend</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:pPr><w:ind w:left=`"720`"/></w:pPr><w:r><w:t>Cascade</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr><w:r><w:t>Flood</w:t></w:r></w:p>"
$outlineBody += "<w:p><w:r><w:t>Foliage &amp; plant representation</w:t></w:r></w:p>"

$outlineTarget = $d.Range($startPos, $endPos)
$outlineTarget.InsertXML((New-PkgXml $outlineBody))

# ---------------------------------------------------------------------------
# 2) Insert the new "Potential advantages/disadvantages of these
#    approaches..." paragraph (split across several runs) plus a trailing
#    blank paragraph, right before "Equations used for fluid dynamics...".
# ---------------------------------------------------------------------------

$rEq = $d.Content
$rEq.Find.Execute("Equations used for fluid dynamics", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$eqStart = $rEq.Start

$advBody = "<w:p>"
$advBody += "<w:r><w:t>Potential advantages/disadvantages of th</w:t></w:r>"
$advBody += "<w:r><w:t>ese</w:t></w:r>"
$advBody += "<w:r><w:t xml:space=`"preserve`"> approach</w:t></w:r>"
$advBody += "<w:r><w:t>es</w:t></w:r>"
$advBody += "<w:r><w:t xml:space=`"preserve`"> compared to traditional methods</w:t></w:r>"
$advBody += "</w:p>"
$advBody += "<w:p/>"

$advTarget = $d.Range($eqStart, $eqStart)
$advTarget.InsertXML((New-PkgXml $advBody))

# ---------------------------------------------------------------------------
# 3) Rename the "Implementation & Results & Analysis" Heading1 to
#    "Critical Reflection".
# ---------------------------------------------------------------------------

$rHead = $d.Content
$rHead.Find.Execute("Implementation & Results & Analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headTarget = $d.Range($rHead.Start, $rHead.End)
$headTarget.InsertXML((New-PkgXml "<w:p><w:r><w:t>Critical Reflection</w:t></w:r></w:p>"))

# ---------------------------------------------------------------------------
# 4) Mark a lastRenderedPageBreak at the start of the "What have I
#    accomplished?" paragraph's run.
# ---------------------------------------------------------------------------

$rAcc = $d.Content
$rAcc.Find.Execute("What have I accomplished? Compare with existing models & simulations in terms of realism & representation. Looking back, would I use a node-based or particle-based simulation?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$accTarget = $d.Range($rAcc.Start, $rAcc.End)
$accBody = "<w:p><w:r><w:lastRenderedPageBreak/><w:t>What have I accomplished? Compare with existing models &amp; simulations in terms of realism &amp; representation. Looking back, would I use a node-based or particle-based simulation?</w:t></w:r></w:p>"
$accTarget.InsertXML((New-PkgXml $accBody))

Write-Output "Edit complete"
